$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Person")

# New header columns
$ws.Range("E1").Value = "Birhdate"
$ws.Range("F1").Value = "Account Value"

# Birhdate column (E) - stored as date serials, formatted as YYYY-MM-DD
$ws.Range("E2").Value = 32874
$ws.Range("E3").Value = 29221
$ws.Range("E4").Value = 27395
$ws.Range("E5").Value = 31048
$ws.Range("E2:E5").NumberFormat = "YYYY\-MM\-DD"

# Account Value column (F) - plain numbers
$ws.Range("F2").Value = 1000
$ws.Range("F3").Value = 2000
$ws.Range("F4").Value = 5000
$ws.Range("F5").Value = 10000

# Update the active selection on the Person sheet to match the target workbook
$ws.Range("J27").Select()
